$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 19.97606611088952
$ws.Range("C2").Value = 8.599995403956191
$ws.Range("D2").Value = 7.737317456828754
$ws.Range("F2").Value = 37.90483260778424
$ws.Range("G2").Value = 43.60082851509443
$ws.Range("H2").Value = 17.9380994333733
$ws.Range("I2").Value = 27.32989560936187
$ws.Range("J2").Value = 10.0497792957749
$ws.Range("L2").Value = 12.24134710680623

$ws.Range("B3").Value = 19.46018987957777
$ws.Range("C3").Value = 8.125619055478154
$ws.Range("D3").Value = 7.732247548079606
$ws.Range("F3").Value = 37.99976338138324
$ws.Range("G3").Value = 43.67005390430137
$ws.Range("H3").Value = 18.00423218823858
$ws.Range("I3").Value = 27.47614163458842
$ws.Range("J3").Value = 10.07572138251402
$ws.Range("L3").Value = 12.22312599818097

$ws.Range("B4").Value = 19.13912887565604
$ws.Range("C4").Value = 7.818335729163691
$ws.Range("D4").Value = 7.729529556533748
$ws.Range("F4").Value = 38.06910172073336
$ws.Range("G4").Value = 43.72919580645645
$ws.Range("H4").Value = 18.0488883786246
$ws.Range("I4").Value = 27.57275885126436
$ws.Range("J4").Value = 10.09248596293991
$ws.Range("L4").Value = 12.21360417817052

$ws.Range("B5").Value = 19.00739900826221
$ws.Range("C5").Value = 7.689144434008319
$ws.Range("D5").Value = 7.728521936733896
$ws.Range("F5").Value = 38.1001243032461
$ws.Range("G5").Value = 43.75745709962156
$ws.Range("H5").Value = 18.06810181854071
$ws.Range("I5").Value = 27.61384196384601
$ws.Range("J5").Value = 10.0995285229267
$ws.Range("L5").Value = 12.21014560647106

$ws.Range("B6").Value = 18.98547703205256
$ws.Range("C6").Value = 7.667454386630306
$ws.Range("D6").Value = 7.728360680311384
$ws.Range("F6").Value = 38.10544230783991
$ws.Range("G6").Value = 43.76240041979576
$ws.Range("H6").Value = 18.07135345652575
$ws.Range("I6").Value = 27.6207669814287
$ws.Range("J6").Value = 10.10071068843276
$ws.Range("L6").Value = 12.20959685445451

$ws.Range("B7").Value = 19.13735568696816
$ws.Range("C7").Value = 7.816609404312063
$ws.Range("D7").Value = 7.72951556170162
$ws.Range("F7").Value = 38.06950891679304
$ws.Range("G7").Value = 43.72956013518591
$ws.Range("H7").Value = 18.04914338916826
$ws.Range("I7").Value = 27.57330599147209
$ws.Range("J7").Value = 10.09258008678052
$ws.Range("L7").Value = 12.21355582401718

$ws.Range("B8").Value = 19.79919706693741
$ws.Range("C8").Value = 8.43978399665391
$ws.Range("D8").Value = 7.735487783885178
$ws.Range("F8").Value = 37.93526420969756
$ws.Range("G8").Value = 43.62123134349037
$ws.Range("H8").Value = 17.9600596582092
$ws.Range("I8").Value = 27.37890268563205
$ws.Range("J8").Value = 10.05855098840917
$ws.Range("L8").Value = 12.23472038138953

$ws.Range("B9").Value = 21.05507946396622
$ws.Range("C9").Value = 9.532893579650446
$ws.Range("D9").Value = 7.750304003829359
$ws.Range("F9").Value = 37.76019318263818
$ws.Range("G9").Value = 43.54170507451144
$ws.Range("H9").Value = 17.81762971428136
$ws.Range("I9").Value = 27.05200884794935
$ws.Range("J9").Value = 9.998424728742711
$ws.Range("L9").Value = 12.28931554659242

$ws.Range("B10").Value = 21.94255358571565
$ws.Range("C10").Value = 10.25557549827767
$ws.Range("D10").Value = 7.763043700401792
$ws.Range("F10").Value = 37.68595509731887
$ws.Range("G10").Value = 43.56532326443926
$ws.Range("H10").Value = 17.73282291273065
$ws.Range("I10").Value = 26.84524436793614
$ws.Range("J10").Value = 9.9582357989598
$ws.Range("L10").Value = 12.33722460013535

$ws.Range("B11").Value = 22.33684032364144
$ws.Range("C11").Value = 10.56664460532499
$ws.Range("D11").Value = 7.769233498175968
$ws.Range("F11").Value = 37.66410498497262
$ws.Range("G11").Value = 43.59403381539821
$ws.Range("H11").Value = 17.69858434924245
$ws.Range("I11").Value = 26.75849802430112
$ws.Range("J11").Value = 9.940809659347149
$ws.Range("L11").Value = 12.3606692737801

$ws.Range("B12").Value = 22.484656631738
$ws.Range("C12").Value = 10.68188297853911
$ws.Range("D12").Value = 7.771633366093374
$ws.Range("F12").Value = 37.6575525910483
$ws.Range("G12").Value = 43.60749608163418
$ws.Range("H12").Value = 17.68624604951725
$ws.Range("I12").Value = 26.72670635173462
$ws.Range("J12").Value = 9.934333275295625
$ws.Range("L12").Value = 12.36978040907299

$ws.Range("B13").Value = 22.45289011309291
$ws.Range("C13").Value = 10.65717822082975
$ws.Range("D13").Value = 7.771114037120245
$ws.Range("F13").Value = 37.65888708310462
$ws.Range("G13").Value = 43.60448147420117
$ws.Range("H13").Value = 17.68887538652334
$ws.Range("I13").Value = 26.73350614004512
$ws.Range("J13").Value = 9.935722640630781
$ws.Range("L13").Value = 12.36780786318566

$ws.Range("B14").Value = 22.34903188682911
$ws.Range("C14").Value = 10.57617664102889
$ws.Range("D14").Value = 7.769429820527947
$ws.Range("F14").Value = 37.6635313699274
$ws.Range("G14").Value = 43.59508942734101
$ws.Range("H14").Value = 17.69755668525623
$ws.Range("I14").Value = 26.75586128086049
$ws.Range("J14").Value = 9.940274391196487
$ws.Range("L14").Value = 12.36141420465694

$ws.Range("B15").Value = 22.28521766401875
$ws.Range("C15").Value = 10.52622753160525
$ws.Range("D15").Value = 7.768405448555174
$ws.Range("F15").Value = 37.6666005577552
$ws.Range("G15").Value = 43.58967398422666
$ws.Range("H15").Value = 17.70295598461147
$ws.Range("I15").Value = 26.76969231503901
$ws.Range("J15").Value = 9.943078408139842
$ws.Range("L15").Value = 12.35752813803634

$ws.Range("B16").Value = 21.91658441205712
$ws.Range("C16").Value = 10.2348892789195
$ws.Range("D16").Value = 7.762647052241434
$ws.Range("F16").Value = 37.68762360781144
$ws.Range("G16").Value = 43.56380914410159
$ws.Range("H16").Value = 17.73514807840267
$ws.Range("I16").Value = 26.85106107655992
$ws.Range("J16").Value = 9.959391813873379
$ws.Range("L16").Value = 12.33572529525927

$ws.Range("B17").Value = 21.68792445316495
$ws.Range("C17").Value = 10.05162230706417
$ws.Range("D17").Value = 7.759214929299107
$ws.Range("F17").Value = 37.70357959049962
$ws.Range("G17").Value = 43.5525494184204
$ws.Range("H17").Value = 17.75601070517299
$ws.Range("I17").Value = 26.90285504949436
$ws.Range("J17").Value = 9.969618392488378
$ws.Range("L17").Value = 12.32276979318955

$ws.Range("B18").Value = 21.55552506395566
$ws.Range("C18").Value = 9.944548141140436
$ws.Range("D18").Value = 7.757278038733677
$ws.Range("F18").Value = 37.71387875007567
$ws.Range("G18").Value = 43.54776404367193
$ws.Range("H18").Value = 17.76841871068196
$ws.Range("I18").Value = 26.93333344319789
$ws.Range("J18").Value = 9.975581055496837
$ws.Range("L18").Value = 12.31547375997825

$ws.Range("B19").Value = 21.51055005097414
$ws.Range("C19").Value = 9.908009403946616
$ws.Range("D19").Value = 7.756628649780341
$ws.Range("F19").Value = 37.71755827323876
$ws.Range("G19").Value = 43.54643395462848
$ws.Range("H19").Value = 17.77268989590717
$ws.Range("I19").Value = 26.94377086943754
$ws.Range("J19").Value = 9.977613772431187
$ws.Range("L19").Value = 12.31303030250584

$ws.Range("B20").Value = 21.71235781607785
$ws.Range("C20").Value = 10.07130370744802
$ws.Range("D20").Value = 7.759576442811618
$ws.Range("F20").Value = 37.70176489879378
$ws.Range("G20").Value = 43.5535729718036
$ws.Range("H20").Value = 17.75374755609394
$ws.Range("I20").Value = 26.89727025851173
$ws.Range("J20").Value = 9.968521417552848
$ws.Range("L20").Value = 12.32413285127813

$ws.Range("B21").Value = 22.37957906635246
$ws.Range("C21").Value = 10.60003824739553
$ws.Range("D21").Value = 7.769923004374793
$ws.Range("F21").Value = 37.66212044795931
$ws.Range("G21").Value = 43.59777776666773
$ws.Range("H21").Value = 17.69498973493246
$ws.Range("I21").Value = 26.74926629260315
$ws.Range("J21").Value = 9.938934111214595
$ws.Range("L21").Value = 12.36328588426457

$ws.Range("B22").Value = 22.80689841203506
$ws.Range("C22").Value = 10.93069065820876
$ws.Range("D22").Value = 7.777010678624226
$ws.Range("F22").Value = 37.64624893588406
$ws.Range("G22").Value = 43.64176637932118
$ws.Range("H22").Value = 17.66024463951168
$ws.Range("I22").Value = 26.65870259255686
$ws.Range("J22").Value = 9.920310953904831
$ws.Range("L22").Value = 12.39023150316655

$ws.Range("B23").Value = 22.57967223755009
$ws.Range("C23").Value = 10.75558266772769
$ws.Range("D23").Value = 7.77319832811679
$ws.Range("F23").Value = 37.65379916940309
$ws.Range("G23").Value = 43.61690607713946
$ws.Range("H23").Value = 17.6784532241916
$ws.Range("I23").Value = 26.70647204798446
$ws.Range("J23").Value = 9.930185354517791
$ws.Range("L23").Value = 12.37572741444322

$ws.Range("B24").Value = 21.70131440931885
$ws.Range("C24").Value = 10.06241107643242
$ws.Range("D24").Value = 7.759412889501085
$ws.Range("F24").Value = 37.7025818144096
$ws.Range("G24").Value = 43.55310496573249
$ws.Range("H24").Value = 17.75476943696456
$ws.Range("I24").Value = 26.89979295817403
$ws.Range("J24").Value = 9.969017100499244
$ws.Range("L24").Value = 12.32351613813727

$ws.Range("B25").Value = 20.72086443563124
$ws.Range("C25").Value = 9.251201496383379
$ws.Range("D25").Value = 7.745967350735469
$ws.Range("F25").Value = 37.79804509810899
$ws.Range("G25").Value = 43.54887219318033
$ws.Range("H25").Value = 17.85268922983849
$ws.Range("I25").Value = 27.13459697123688
$ws.Range("J25").Value = 10.01398760648383
$ws.Range("L25").Value = 12.27316251582346
